$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the style used by the other
# header cells (e.g. G1: bold, bordered, centered/top-aligned).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# Numeric value for the new column on the data row.
$ws.Range("H2").Value = 1
